$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = 0.7790697674418605
$ws.Cells.Item(3, 3).Value = 0.813953488372093
$ws.Cells.Item(3, 4).Value = 0.872093023255814
$ws.Cells.Item(3, 5).Value = 0.872093023255814

$ws.Cells.Item(4, 2).Value = 0.7294117647058823
$ws.Cells.Item(4, 3).Value = 0.7411764705882353
$ws.Cells.Item(4, 4).Value = 0.8235294117647058
$ws.Cells.Item(4, 5).Value = 0.8117647058823529

$ws.Cells.Item(5, 2).Value = 0.7058823529411765
$ws.Cells.Item(5, 3).Value = 0.7176470588235294
$ws.Cells.Item(5, 4).Value = 0.7764705882352941
$ws.Cells.Item(5, 5).Value = 0.8117647058823529

$ws.Cells.Item(6, 2).Value = 0.7647058823529411
$ws.Cells.Item(6, 3).Value = 0.8
$ws.Cells.Item(6, 4).Value = 0.8705882352941177
$ws.Cells.Item(6, 5).Value = 0.8705882352941177

$ws.Cells.Item(7, 2).Value = 0.8
$ws.Cells.Item(7, 3).Value = 0.7764705882352941
$ws.Cells.Item(7, 4).Value = 0.8705882352941177
$ws.Cells.Item(7, 5).Value = 0.8823529411764706

$ws.Cells.Item(8, 2).Value = 0.8604651162790697
$ws.Cells.Item(8, 3).Value = 0.8604651162790697
$ws.Cells.Item(8, 4).Value = 0.8953488372093024
$ws.Cells.Item(8, 5).Value = 0.8953488372093024

$ws.Cells.Item(9, 2).Value = 0.7790697674418605
$ws.Cells.Item(9, 3).Value = 0.8372093023255814
$ws.Cells.Item(9, 4).Value = 0.872093023255814
$ws.Cells.Item(9, 5).Value = 0.8837209302325582

$ws.Cells.Item(10, 2).Value = 0.7906976744186046
$ws.Cells.Item(10, 3).Value = 0.8023255813953488
$ws.Cells.Item(10, 4).Value = 0.8488372093023255
$ws.Cells.Item(10, 5).Value = 0.8604651162790697

$ws.Cells.Item(11, 2).Value = 0.7674418604651163
$ws.Cells.Item(11, 3).Value = 0.7790697674418605
$ws.Cells.Item(11, 4).Value = 0.8255813953488372
$ws.Cells.Item(11, 5).Value = 0.813953488372093

$ws.Cells.Item(12, 2).Value = 0.813953488372093
$ws.Cells.Item(12, 3).Value = 0.8255813953488372
$ws.Cells.Item(12, 4).Value = 0.8604651162790697
$ws.Cells.Item(12, 5).Value = 0.8372093023255814

[void]$ws.Range("E3:E12").Select()
